# pH_UL.xlsx update: rename the "Run" header to "Model_", add the *_new
# metric columns (T:AK) plus the t_value_pH_UL column, and append two new
# rows (Model2 / Model3) with their corresponding metric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
$ws.Range("A1").Value = "Model_"

$newHeaders = @(
    "RMSE_X_new", "NMRSE_X_new", "MAPE_X_new",
    "RMSE_C_new", "NMRSE_C_new", "MAPE_C_new",
    "RMSE_N_new", "NMRSE_N_new", "MAPE_N_new",
    "RMSE_pH_new", "NMRSE_pH_new", "MAPE_pH_new",
    "AIC_new", "BIC_new",
    "RMSE_new", "NMRSE_new", "MAPE_new",
    "t_value_pH_UL"
)

# Columns T..AK are columns 20..37
for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $ws.Cells.Item(1, 20 + $i).Value = $newHeaders[$i]
}

# Match the look of the existing header cells (bold, centered, boxed) by
# cloning the formatting that's already on the row-1 header cells.
$ws.Range("A1").Copy()
$headerRange = $ws.Range($ws.Cells.Item(1, 20), $ws.Cells.Item(1, 37))
$headerRange.PasteSpecial(-4122)

# --- Row 3: Model2 ------------------------------------------------------
$ws.Range("A3").Value = "Model2"
$ws.Range("B3").Value = 6.999999490197824

$row3 = @(
    0.1916045932598129, 0.1953173813621255, 15.73112756341508,
    0.4996824741820515, 0.2193064667095056, 9.109455022142289,
    0.07267687114484886, 0.572933716762807, 7.324738627686737,
    0.1538708269774571, 0.1947731987056419, 1.784108967285242,
    -89.45284559428846, -87.86932665583235,
    0.2807815227882198, 0.04034032195096889, 8.487357545132337,
    72116.69311008597
)
for ($i = 0; $i -lt $row3.Length; $i++) {
    $ws.Cells.Item(3, 20 + $i).Value = $row3[$i]
}

# --- Row 4: Model3 ------------------------------------------------------
$ws.Range("A4").Value = "Model3"
$ws.Range("B4").Value = 6.999997130806487

$row4 = @(
    0.1916046854664078, 0.1953174753554393, 15.73114390224424,
    0.4996830423199394, 0.219306716060482, 9.109464646988483,
    0.0726768949811793, 0.5729339046717962, 7.324740893334438,
    0.1538708788262824, 0.1947732643370664, 1.784110473202808,
    -89.45277452742999, -87.86925558897389,
    0.2807817999308668, 0.04034036176848785, 8.487364978942493,
    72117.88406508609
)
for ($i = 0; $i -lt $row4.Length; $i++) {
    $ws.Cells.Item(4, 20 + $i).Value = $row4[$i]
}

# --- Blank placeholder cells -------------------------------------------
# In the source workbook, row 2 also carries empty string placeholders in
# T2:AK2, and rows 3/4 carry empty string placeholders in C:S (the columns
# that row 2/the "og" metrics occupy but rows 3/4 do not use). Write empty
# strings so these cells are touched/created where the host supports it.
for ($c = 20; $c -le 37; $c++) {
    $ws.Cells.Item(2, $c).Value = ""
}
for ($c = 3; $c -le 19; $c++) {
    $ws.Cells.Item(3, $c).Value = ""
    $ws.Cells.Item(4, $c).Value = ""
}
